$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching style of existing header (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-61: row, I value, J value
$data = @(
    @(2,7,7),
    @(3,7,7),
    @(4,3,4),
    @(5,8,8),
    @(6,9,9),
    @(7,7,7),
    @(8,1,1),
    @(9,10,10),
    @(10,5,5),
    @(11,4,5),
    @(12,8,8),
    @(13,6,6),
    @(14,7,7),
    @(15,6,6),
    @(16,7,7),
    @(17,6,6),
    @(18,8,8),
    @(19,5,5),
    @(20,4,6),
    @(21,9,9),
    @(22,7,7),
    @(23,5,6),
    @(24,8,8),
    @(25,4,5),
    @(26,8,8),
    @(27,5,6),
    @(28,7,7),
    @(29,7,7),
    @(30,6,6),
    @(31,10,11),
    @(32,5,5),
    @(33,6,6),
    @(34,6,7),
    @(35,10,10),
    @(36,9,9),
    @(37,5,7),
    @(38,7,7),
    @(39,8,8),
    @(40,9,9),
    @(41,8,8),
    @(42,8,8),
    @(43,5,6),
    @(44,9,9),
    @(45,5,5),
    @(46,7,8),
    @(47,7,7),
    @(48,7,7),
    @(49,7,7),
    @(50,6,6),
    @(51,7,7),
    @(52,7,7),
    @(53,7,8),
    @(54,6,6),
    @(55,7,7),
    @(56,8,8),
    @(57,8,8),
    @(58,7,7),
    @(59,9,9),
    @(60,8,8),
    @(61,8,8)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
